# "Generate Report for Archive"
# Status moved from "Ready for handoff" to "In Translation" on the single
# tracked file, across the Overview roll-up sheet and each per-locale
# detail sheet. Excel then re-fit the (now shorter) Status columns.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale "Status" columns (E = zh-cn, F = de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn detail sheet: "Status" column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: "Status" column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
